$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.912.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.499.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.584'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("E9").Value = '  +4.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.430'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.105.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.930.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.489.67'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '392.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.534'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.180'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.996'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.13'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.05%  '
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.66'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.36'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.61'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '164.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.877'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.39%  '
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.31'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.92%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.846.70'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.22%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0735'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.77%  '
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.53'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0302'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '339.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.07'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.840'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.44%  '
